$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Copy the formatting of the last existing data row (126) down into
#        the seven new rows (127-133) so columns A and C:H inherit the same
#        number formats / fonts as the rest of the log. ---
$ws.Range("A126:J126").Copy()
$ws.Range("A127:J133").PasteSpecial(-4122)

# --- 2. New daily entries for 07/05/2018 - 13/05/2018 ---
$data = @(
    @(43227,1522,43,20,143,2691,136,3250),
    @(43228,1568,42,30,166,2811,137,4000),
    @(43229,1539,64,17,131,2694,114,5000),
    @(43230,1507,35,25,180,2938,118,4500),
    @(43231,1473,43,15,173,2805,75,6500),
    @(43232,1806,88,16,142,3333,110,3000),
    @(43233,1765,41,17,240,3214,97,1500)
)

$r = 127
foreach ($row in $data) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $ws.Cells.Item($r,7).Value = $row[6]
    $ws.Cells.Item($r,8).Value = $row[7]
    $ws.Cells.Item($r,9).Formula = '=IF(H' + $r + '>=2200,"Yes","No")'
    $ws.Cells.Item($r,10).Formula = '=IF(B' + $r + '<=1800,"Yes","No")'
    $r = $r + 1
}

# --- 3. Column B on the new rows uses a dedicated "0 w/ padding" number
#        format (matches a Comma-style cell whose decimals were reduced to 0). ---
$ws.Range("B127:B133").NumberFormatLocal = "0_ ;\-0\ "

# --- 4. Update the view: scrolled down so the newly-added rows are visible,
#        with the final entry selected. ---
$win = $excel.ActiveWindow
$win.ScrollRow = 121
$win.ScrollColumn = 1
$ws.Range("B133").Select()
